# Applies the "Added handling of common packages." commit:
# The classFields sheet rows were re-ordered (within each class's
# field-group) to reflect the new field declaration order picked up
# by the structure-mining tool. Column A (Class Name) and column C
# (Field Modifier) never change per group; only column B (Field Name)
# and column D (Field Type) are permuted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# com.macro.mall.auth.constant.MessageConstant fields (rows 2-8)
$ws.Cells.Item(3, 2).Value = "PERMISSION_DENIED"
$ws.Cells.Item(4, 2).Value = "ACCOUNT_DISABLED"
$ws.Cells.Item(5, 2).Value = "USERNAME_PASSWORD_ERROR"
$ws.Cells.Item(6, 2).Value = "ACCOUNT_LOCKED"
$ws.Cells.Item(7, 2).Value = "CREDENTIALS_EXPIRED"
$ws.Cells.Item(8, 2).Value = "ACCOUNT_EXPIRED"

# com.macro.mall.auth.config.Oauth2ServerConfig fields (rows 10-13)
$ws.Cells.Item(10, 2).Value = "jwtTokenEnhancer"
$ws.Cells.Item(10, 4).Value = "com.macro.mall.auth.component.JwtTokenEnhancer"
$ws.Cells.Item(13, 2).Value = "userDetailsService"
$ws.Cells.Item(13, 4).Value = "com.macro.mall.auth.service.impl.UserServiceImpl"

# com.macro.mall.auth.domain.Oauth2TokenDto fields (rows 14-17)
$ws.Cells.Item(14, 2).Value = "tokenHead"
$ws.Cells.Item(14, 4).Value = "java.lang.String"
$ws.Cells.Item(17, 2).Value = "expiresIn"
$ws.Cells.Item(17, 4).Value = "int"

# com.macro.mall.auth.service.impl.UserServiceImpl fields (rows 18-20)
$ws.Cells.Item(19, 2).Value = "request"
$ws.Cells.Item(19, 4).Value = "javax.servlet.http.HttpServletRequest"
$ws.Cells.Item(20, 2).Value = "memberService"
$ws.Cells.Item(20, 4).Value = "com.macro.mall.auth.service.UmsMemberService"

# com.macro.mall.auth.domain.Oauth2TokenDto$Oauth2TokenDtoBuilder fields (rows 21-24)
$ws.Cells.Item(21, 2).Value = "token"
$ws.Cells.Item(22, 2).Value = "expiresIn"
$ws.Cells.Item(22, 4).Value = "int"
$ws.Cells.Item(23, 2).Value = "refreshToken"
$ws.Cells.Item(24, 2).Value = "tokenHead"
$ws.Cells.Item(24, 4).Value = "java.lang.String"

# com.macro.mall.auth.domain.SecurityUser fields (rows 25-30)
$ws.Cells.Item(25, 2).Value = "authorities"
$ws.Cells.Item(25, 4).Value = "java.util.Collection"
$ws.Cells.Item(26, 2).Value = "id"
$ws.Cells.Item(26, 4).Value = "java.lang.Long"
$ws.Cells.Item(27, 2).Value = "clientId"
$ws.Cells.Item(28, 2).Value = "username"
$ws.Cells.Item(30, 2).Value = "password"
